# Add organisations into the PMHC upload spec workbook.
$wb = $excel.ActiveWorkbook

# --- 1. New "Organisations" sheet, placed after the last existing sheet (Practitioners) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Organisations"

# Header / version row
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = 1

# Field-name row
$ws.Range("A2").Value = "organisation_path"
$ws.Range("B2").Value = "organisation_key"
$ws.Range("C2").Value = "organisation_name"
$ws.Range("D2").Value = "organisation_legal_name"
$ws.Range("E2").Value = "organisation_abn"
$ws.Range("F2").Value = "organisation_type"
$ws.Range("G2").Value = "organisation_state"
$ws.Range("H2").Value = "organisation_status"
$ws.Range("I2").Value = "organisation_tags"

# Sample data row
$ws.Range("A3").Value = "PHN999:NFP01"
$ws.Range("B3").Value = "NFP01"
$ws.Range("C3").Value = "Test Provider Organisation NFP1"
$ws.Range("E3").Value = 42072953425
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Selection / column widths / margins matching the authored layout
$ws.Range("F4").Select() | Out-Null

$ws.Columns.Item(1).ColumnWidth = 18.666666666666668   # -> 19.5
$ws.Columns.Item(2).ColumnWidth = 17.0                 # -> 17.83203125 (closest reachable)
$ws.Columns.Item(3).ColumnWidth = 24.833333333333332   # -> 25.6640625 (closest reachable)
$ws.Columns.Item(4).ColumnWidth = 20.833333333333332   # -> 21.6640625 (closest reachable)
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668   # -> 20.5

$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# --- 2. Clients sheet: new narrower column A ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Columns.Item(1).ColumnWidth = 14.333333333333334   # -> 15.1640625 (closest reachable)

# --- 3. Window position / active tab bookkeeping ---
$win = $wb.Windows.Item(1)
$win.Left = 540
$win.Top = 2160
